$d = $word.ActiveDocument

$replacements = @(
    @("913×6=", "954×8="),
    @("309×5=", "394×9="),
    @("714×7=", "987×6="),
    @("425×9=", "486×2="),
    @("780×7=", "968×2="),
    @("786×4=", "726×4="),
    @("429×3=", "315×5="),
    @("902×8=", "599×4="),
    @("219×9=", "488×8="),
    @("209×4=", "737×4="),
    @("237×2=", "346×3="),
    @("834×7=", "498×7="),
    @("868×5=", "164×7="),
    @("250×8=", "164×4="),
    @("608×5=", "230×6="),
    @("653×2=", "329×3="),
    @("259×6=", "334×2="),
    @("196×3=", "771×4="),
    @("244×3=", "479×7="),
    @("686×9=", "517×6="),
    @("303×2=", "630×4="),
    @("112×2=", "787×4="),
    @("912×7=", "648×7="),
    @("388×6=", "286×9="),
    @("652×5=", "354×4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
